$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 12.66339259259259
$ws.Range("F2").Value = 1.725244444444445
$ws.Range("G2").Value = -0.03344951089196524
$ws.Range("E3").Value = 12.66056296296296
$ws.Range("F3").Value = 1.725111111111111
$ws.Range("G3").Value = -0.03321858704313385
$ws.Range("E4").Value = 12.6610962962963
$ws.Range("F4").Value = 1.725377777777778
$ws.Range("G4").Value = -0.03326211195704976
$ws.Range("E5").Value = 13.05798737384537
$ws.Range("F5").Value = 1.705089338375171
$ws.Range("G5").Value = -0.06565208067762773
$ws.Range("E6").Value = 13.05530945381646
$ws.Range("F6").Value = 1.704919375040773
$ws.Range("G6").Value = -0.06543353773000549
$ws.Range("E7").Value = 13.05530945381646
$ws.Range("F7").Value = 1.704919375040773
$ws.Range("G7").Value = -0.06543353773000549
$ws.Range("E8").Value = 11.51160266666667
$ws.Range("F8").Value = 1.70644
$ws.Range("G8").Value = 0.06054715918637621
$ws.Range("E9").Value = 11.51160266666667
$ws.Range("F9").Value = 1.706573333333333
$ws.Range("G9").Value = 0.06054715918637621
$ws.Range("E10").Value = 11.51148266666667
$ws.Range("F10").Value = 1.706613333333333
$ws.Range("G10").Value = 0.06055695229200742
$ws.Range("E11").Value = 12.01697927585532
$ws.Range("F11").Value = 1.686726009705235
$ws.Range("G11").Value = 0.01930377154255025
$ws.Range("E12").Value = 12.01688363015876
$ws.Range("F12").Value = 1.686736041782117
$ws.Range("G12").Value = 0.01931157711262943
$ws.Range("E13").Value = 12.0165932839573
$ws.Range("F13").Value = 1.686701250062081
$ws.Range("G13").Value = 0.01933527203780017
$ws.Range("E14").Value = 10.807584
$ws.Range("F14").Value = 1.68992
$ws.Range("G14").Value = 0.1180015689273384
$ws.Range("E15").Value = 10.805792
$ws.Range("F15").Value = 1.689786666666667
$ws.Range("G15").Value = 0.1181478126380958
$ws.Range("E16").Value = 10.80500133333334
$ws.Range("F16").Value = 1.689853333333333
$ws.Range("G16").Value = 0.1182123383229761
$ws.Range("E17").Value = 11.23820993432772
$ws.Range("F17").Value = 1.66213077382702
$ws.Range("G17").Value = 0.08285852507440605
$ws.Range("E18").Value = 11.23719499019746
$ws.Range("F18").Value = 1.661961733207851
$ws.Range("G18").Value = 0.08294135386671719
$ws.Range("E19").Value = 11.23648262536601
$ws.Range("F19").Value = 1.661990930322566
$ws.Range("G19").Value = 0.08299948940040203
$ws.Range("E20").Value = 10.55505983333333
$ws.Range("F20").Value = 1.69519
$ws.Range("G20").Value = 0.1386098675820544
$ws.Range("E21").Value = 10.55576883333333
$ws.Range("F21").Value = 1.69523
$ws.Range("G21").Value = 0.1385520066496176
$ws.Range("E22").Value = 10.556425
$ws.Range("F22").Value = 1.695233333333333
$ws.Range("G22").Value = 0.1384984574039655
$ws.Range("E23").Value = 10.87598456480181
$ws.Range("F23").Value = 1.649513239720326
$ws.Range("G23").Value = 0.1124194526245936
$ws.Range("E24").Value = 10.87659164397954
$ws.Range("F24").Value = 1.6495658273368
$ws.Range("G24").Value = 0.1123699093703109
$ws.Range("E25").Value = 10.87640283862705
$ws.Range("F25").Value = 1.649517985825406
$ws.Range("G25").Value = 0.1123853176266498
$ws.Range("E26").Value = 10.55603874666667
$ws.Range("F26").Value = 1.726937333333333
$ws.Range("G26").Value = 0.1385299792347461
$ws.Range("E27").Value = 10.55567253333333
$ws.Range("F27").Value = 1.726928
$ws.Range("G27").Value = 0.1385598656168866
$ws.Range("E28").Value = 10.5555852
$ws.Range("F28").Value = 1.726884
$ws.Range("G28").Value = 0.1385669928215403
$ws.Range("E29").Value = 10.74822670683174
$ws.Range("F29").Value = 1.655656591032269
$ws.Range("G29").Value = 0.1228456709437671
$ws.Range("E30").Value = 10.74794325744123
$ws.Range("F30").Value = 1.655648019551619
$ws.Range("G30").Value = 0.1228688030256203
$ws.Range("E31").Value = 10.74797729099025
$ws.Range("F31").Value = 1.65564050323588
$ws.Range("G31").Value = 0.1228660255744491
